$d = $word.ActiveDocument

# Locate the "Changes made" text so we can insert the new paragraph right
# after it (and before the trailing _GoBack bookmark that currently sits
# at the end of that same paragraph).
$findRange = $d.Content
$found = $findRange.Find.Execute("Changes made")
$endOfChangesMade = $findRange.End

# Insert a new paragraph break right after "Changes made" (i.e. before the
# _GoBack bookmark, which stays collapsed at that same point and will end
# up in the newly created paragraph).
$breakPoint = $d.Range($endOfChangesMade, $endOfChangesMade)
$breakPoint.InsertParagraphAfter()

# Build up "3" + "rd" (superscript) + " commit" in the freshly created
# paragraph, run by run, so the formatting boundaries match the target.
$cursor = $d.Range($endOfChangesMade + 1, $endOfChangesMade + 1)
$cursor.InsertAfter("3")

$rdRun = $d.Range($cursor.End, $cursor.End)
$rdRun.InsertAfter("rd")
$rdRun.Font.Superscript = $true

$tailRun = $d.Range($rdRun.End, $rdRun.End)
$tailRun.InsertAfter(" commit")
$endOfNewText = $tailRun.End

# Move the "_GoBack" bookmark so it ends up collapsed right after " commit"
# (matching where Word leaves it once you've typed past the old location).
# A collapsed Range placed exactly at a paragraph's last character position
# confuses bookmark placement, so insert a throwaway character after the
# target point first, anchor the bookmark there, then remove the throwaway
# character again - the bookmark stays put.
$placeholder = $d.Range($endOfNewText, $endOfNewText)
$placeholder.InsertAfter("X")

$bookmarkSpot = $d.Range($endOfNewText, $endOfNewText)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$placeholderChar = $d.Range($endOfNewText, $endOfNewText + 1)
$placeholderChar.Delete()
